# Actualización automática 2025-09-08 09:55:08
# Refresh the "CUMPLIMIENTO MENSUAL" sheet data: 4 discontinued product
# groups (GRANITO, LED, PANELES PU, PANELES PVC) are removed and every
# remaining row is refreshed with new PRESUPUESTO / VENTA / POR CUMPLIR /
# CUMPLIMIENTO figures, including a recomputed TOTAL row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# 1) Drop the rows for groups that no longer appear in the refreshed pull.
#    Delete from the bottom up so earlier row numbers stay valid.
$ws.Rows.Item(14).Delete()   # PANELES PVC
$ws.Rows.Item(13).Delete()   # PANELES PU
$ws.Rows.Item(9).Delete()    # LED
$ws.Rows.Item(5).Delete()    # GRANITO

# 2) Narrow/widen the PRESUPUESTO/POR CUMPLIR/CUMPLIMIENTO columns to
#    their new widths (stored width = ColumnWidth + 5/6).
$ws.Columns.Item(4).ColumnWidth = 13.166666666666666   # -> stored 14
$ws.Columns.Item(5).ColumnWidth = 22.166666666666668   # -> stored 23
$ws.Columns.Item(6).ColumnWidth = 24.166666666666668   # -> stored 25

# 3) Refresh every remaining data row (labels already line up after the
#    deletions above) with the new figures.
$ws.Range("B2").Value = "240X120 PORCELANATO"
$ws.Range("C2").Value = 743.710083532391
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 743.710083532391
$ws.Range("F2").Value = 0

$ws.Range("B3").Value = "240X80 PORCELANATO"
$ws.Range("C3").Value = 5504.61890386263
$ws.Range("D3").Value = 475.2
$ws.Range("E3").Value = 5029.41890386263
$ws.Range("F3").Value = 0.08632750210310632

$ws.Range("B4").Value = "FREGADEROS DE COCINA"
$ws.Range("C4").Value = 521.61144263264
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 521.61144263264
$ws.Range("F4").Value = 0

$ws.Range("B5").Value = "GRIFERIAS"
$ws.Range("C5").Value = 150
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 150
$ws.Range("F5").Value = 0

$ws.Range("B6").Value = "INODOROS"
$ws.Range("C6").Value = 2907.58368146026
$ws.Range("D6").Value = 1803.6
$ws.Range("E6").Value = 1103.98368146026
$ws.Range("F6").Value = 0.6203088879265507

$ws.Range("B7").Value = "LAVABOS"
$ws.Range("C7").Value = 886.711016287574
$ws.Range("D7").Value = 127.8
$ws.Range("E7").Value = 758.9110162875741
$ws.Range("F7").Value = 0.1441281292918464

$ws.Range("B8").Value = "NO RESURTIBLES"
$ws.Range("C8").Value = 516.121873547834
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 516.121873547834
$ws.Range("F8").Value = 0

$ws.Range("B9").Value = "OTROS"
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 0

$ws.Range("B10").Value = "PANELES DECORATIVOS"
$ws.Range("C10").Value = 388.107983534392
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 388.107983534392
$ws.Range("F10").Value = 0

$ws.Range("B11").Value = "PIEDRA SINTERIZADA"
$ws.Range("C11").Value = 5844.44916370549
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 5844.44916370549
$ws.Range("F11").Value = 0

$ws.Range("B12").Value = "PORCELANATO"
$ws.Range("C12").Value = 36823.6430921171
$ws.Range("D12").Value = 21833.23
$ws.Range("E12").Value = 14990.4130921171
$ws.Range("F12").Value = 0.5929133612712502

$ws.Range("B13").Value = "PUERTAS DE SEGURIDAD"
$ws.Range("C13").Value = 222.087330240682
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 222.087330240682
$ws.Range("F13").Value = 0

$ws.Range("B14").Value = "SAL SOLUBLE"
$ws.Range("C14").Value = 916.0985952127839
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 916.0985952127839
$ws.Range("F14").Value = 0

# 4) TOTAL row.
$ws.Range("B15").Value = "TOTAL"
$ws.Range("C15").Value = 55424.74316613378
$ws.Range("D15").Value = 24239.83
$ws.Range("E15").Value = 31184.91316613378
$ws.Range("F15").Value = 0.437346726665777
